$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 38, which shifts the old row 38
# ("Extra Tall Female Stackable Header ...") down to row 39, preserving its
# formatting and values, and leaves a new blank row 38 (the separator).
$ws.Rows("38").Insert()

# Populate the now-available blank row 37 (previously an empty separator row)
# with the new "696ZZ Bearing" line item, copying the style from row 36 so the
# number/text formatting matches the rest of the table.
$ws.Range("A36:C36").Copy()
$ws.Range("A37:C37").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A37").Value = 2
$ws.Range("B37").Value = "696ZZ Bearing 6mm x 15mm x 5mm"
$ws.Range("C37").Value = "https://www.amazon.com/gp/product/B07FW26HD4/"

# Match the final selection state recorded in the workbook.
$ws.Range("A37:C37").Select()
